# TMTC0032668 - Changed Test Data for LV Activities - 16 Sep 2024
#
# The Contact sheet's sample "external contact" row is updated to use
# activity-specific test data instead of the generic standard test data,
# and the Contact sheet becomes the active tab (previously UpdateActivity
# was active).

$wb = $excel.ActiveWorkbook

$contact = $wb.Worksheets.Item("Contact")

# Update the sample contact/company values used by this row of test data.
$contact.Range("A2").Value = "Activity Test External Contact"
$contact.Range("B2").Value = "ActivityCompany"

# Make the Contact sheet the active sheet/tab, with A2:B2 selected
# (activeCell A2), matching the new authored state of the workbook.
$contact.Activate()
$contact.Range("A2:B2").Select()
